# Apply updated cryptos list values to worksheet (cryptos.xlsx)
# Commit: Updated cryptos list on Tue Sep 19 06:38:23 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values (column D) are plain decimal numbers
# (e.g. '2.10', '0.0958'). Left alone, Excel's COM layer would coerce a
# plain numeric-looking string into a Double and drop significant trailing
# zeros / formatting (e.g. '2.10' -> 2.1). Force those specific cells to
# Text format first so the literal string is preserved, matching the
# original inlineStr text content.
$textCells = @("D5", "D6", "D10", "D14", "D16", "D19", "D21", "D22", "D23", "D24", "D25", "D27", "D30", "D34", "D38", "D39", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.839.92'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.637.71'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("D5").Value = '216.88'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").Value = '0.508'
$ws.Range("E6").Value = '  +2.10%  '
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").Value = '19.85'
$ws.Range("E10").Value = '  +3.00%  '
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '1.866.57'
$ws.Range("D13").Value = '1.643.57'
$ws.Range("E13").Value = '  +0.69%  '
$ws.Range("D14").Value = '4.12'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '67.16'
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("D17").Value = '26.831.96'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = '217.91'
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("D21").Value = '6.77'
$ws.Range("E21").Value = '  +1.29%  '
$ws.Range("D22").Value = '4.39'
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").Value = '2.43'
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").Value = '9.14'
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("D25").Value = '146.95'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  -0.49%  '
$ws.Range("D27").Value = '7.31'
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '0.0501'
$ws.Range("E30").Value = '  -1.28%  '
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("E33").Value = '  -0.43%  '
$ws.Range("D34").Value = '1.57'
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").Value = '1.266.87'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  +1.72%  '
$ws.Range("D38").Value = '0.533'
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = '0.833'
$ws.Range("E39").Value = '  +1.78%  '
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("D43").Value = '1.778.39'
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '61.75'
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = '2.10'
$ws.Range("E45").Value = '  -0.61%  '
$ws.Range("D46").Value = '91.64'
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("E47").Value = '  -1.08%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0512'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.64'
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '0.0958'
$ws.Range("E51").Value = '  -0.90%  '
